$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row 8: fantasy-team names above the existing player-name row 9 ---
# (row 8 did not exist before; this does NOT shift any existing rows)
$ws.Range("D8").Value = "DAMU DREAM TEAM"
$ws.Range("G8").Value = "JUSTIN CHALLENGERS"
$ws.Range("J8").Value = "Devilish 11"
$ws.Range("M8").Value = "RENGAN25QR"
$ws.Range("P8").Value = "Sundar Night Fury"
$ws.Range("S8").Value = "speedsterse7en"

# Merge each name over its score column, matching the row-9 header layout
$ws.Range("D8:E8").Merge()
$ws.Range("G8:H8").Merge()
$ws.Range("J8:K8").Merge()
$ws.Range("M8:N8").Merge()
$ws.Range("P8:Q8").Merge()
$ws.Range("S8:T8").Merge()

# Style row 8 like the "20% - Accent6" theme, bold, centered, bordered
$row8Rng = $ws.Range("D8:E8,G8:H8,J8:K8,M8:N8,P8:Q8,S8:T8")
$row8Rng.Style = "20% - Accent6"
$row8Rng.Font.Bold = $true
$row8Rng.HorizontalAlignment = -4108
$row8Rng.Borders.LineStyle = 1
$row8Rng.Borders.Weight = 2

# --- Contest 2: DC vs KXI (row 11) results entry ---
$ws.Range("E11").Value = 0
$ws.Range("H11").Value = 60
$ws.Range("K11").Value = 20
$ws.Range("N11").Value = 80
$ws.Range("Q11").Value = 40
$ws.Range("T11").Value = 100

# Selection as left by the editor
$ws.Range("N28").Select()
